# daily auto push: 2026-02-11 19:18 UTC
# A new daily-log row for 2026/02/12 (Thu) was inserted into Sheet1 right
# after the existing 2026/02/11 entries (row 808), pushing every row from
# the old row 808 onward down by one (old last row 849 -> 850).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 808..849 down to 809..850 and open up a blank row 808.
$ws.Rows(808).Insert()

# Fill the new row. Column A holds literal "yyyy/mm/dd" text (not a real
# date) in this sheet, so force text entry (leading apostrophe) and then
# reset the style back to the sheet's default "Normal" so no stray
# number-format style gets attached to the cell.
$ws.Range("A808").Value = "'2026/02/12"
$ws.Range("A808").Style = "Normal"
$ws.Range("B808").Value = "木"
$ws.Range("C808").Value = 2
$ws.Range("D808").Value = 201
